# Update the "last autoupdated" date placeholder text (shown as a
# datetimeFigureOut field) on the slide master and every slide layout,
# and update a label on slide 1 from "GAE " to "Google Cloud ".

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($sh.Type -eq -2147483648 -or $sh.Type -eq 14) {
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            } catch {
            }
        }
        if ($isDatePlaceholder) {
            # Whole-range assign via TextFrame2 - the date placeholder's
            # paragraph is a single auto-updating field occupying the
            # entire range, so this replaces it in place (and mirrors
            # the smtClean="0" bookkeeping PowerPoint itself adds when
            # it re-writes a placeholder run).
            $sh.TextFrame2.TextRange.Text = $newText
        }
    }
}

$newDate = "8/14/2020"

# Slide master.
$master = $p.Designs.Item(1).SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

# Every custom (slide) layout under the master.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# Slide 1: "GAE " -> "Google Cloud " in the "Flowchart: Direct Access
# Storage 84" shape (first run of its text body).
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.Name -eq "Flowchart: Direct Access Storage 84") {
        $tr = $sh.TextFrame.TextRange
        $run1 = $tr.Characters(1, 4)
        if ($run1.Text -eq "GAE ") {
            $run1.Text = "Google Cloud "
        }
    }
}
